$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: add formula to E2 (replacing the "insert Excel formula here" placeholder text) ---
$ws.Range("E2").Formula = "=(B2/C2)*D2"

# --- Row 3: fill in miles/MPG/price inputs ---
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 30
$ws.Range("D3").Value = 1

# --- Row 4: fill in miles/MPG/price inputs ---
$ws.Range("B4").Value = 40
$ws.Range("C4").Value = 30
$ws.Range("D4").Value = 1

# --- Row 5: fill in miles/MPG/price inputs ---
$ws.Range("B5").Value = 29
$ws.Range("C5").Value = 30
$ws.Range("D5").Value = 1

# --- New rows 6-8: additional test cases (values only, styling comes below) ---
$ws.Range("A6").Value = "20 miles left to go, price is 2"
$ws.Range("A7").Value = "Price is 5, we need 2 gallons"
$ws.Range("A8").Value = "5 miles left to go, price is 5"

# Apply the normal data-row style (same style already used for B3:D5) onto the new B6:D8 cells
$ws.Range("B3:D3").Copy()
$ws.Range("B6:D8").PasteSpecial(-4122)
$ws.Range("B6").Value = 20
$ws.Range("C6").Value = 30
$ws.Range("D6").Value = 2
$ws.Range("B7").Value = 60
$ws.Range("C7").Value = 30
$ws.Range("D7").Value = 5
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 30
$ws.Range("D8").Value = 5

# Build the new "highlight" look (base data font + yellow fill + centered) once on a scratch
# cell, then stamp it onto A6:A8 via copy/paste so the style table only gains one new entry.
$scratch = $ws.Range("Z1")
$ws.Range("A2").Copy()
$scratch.PasteSpecial(-4122)
$scratch.Interior.Color = 65535
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4108

$scratch.Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = "20 miles left to go, price is 2"

$scratch.Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = "Price is 5, we need 2 gallons"

$scratch.Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "5 miles left to go, price is 5"

$scratch.Clear()

# --- Shared formula for E3:E8 (out:cost) ---
$ws.Range("E3:E8").Formula = "=(B3/C3)*D3"

# --- Expand the table to cover the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E8"))
